$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the google ad-click URL in B2
$ws.Range("B2").Value = "https://www.google.com/aclk?sa=l&ai=DChcSEwjhw_eAgP2CAxWVXpEFHfP2DxEYABAAGgJscg&ase=2&gclid=EAIaIQobChMI4cP3gID9ggMVlV6RBR3z9g8REBAYASAAEgLelPD_BwE&sig=AOD64_3N2VaRGe2qXXt-6MBEf3xs6xeQMw&adurl&ctype=99"

# Insert a new row at row 4 for "sidedishmedia", shifting the rest down
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = "sidedishmedia"
$ws.Range("B4").Value = "https://www.sidedishmedia.co.uk/"
$ws.Range("C4").Value = "London"
$ws.Range("D4").Value = "UK"
$ws.Range("E4").Value = "hello@sidedishmedia.co.uk"
